$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B ("ID Competição") rows 2-71 were dropped/truncated to 39
# and need to be recovered/corrected to 239.
$ws.Range("B2:B71").Value = 239
